$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Capture original row 2 values before we overwrite anything
$url = $ws.Range("A2").Value2
$user = $ws.Range("B2").Value2
$pass = $ws.Range("C2").Value2
$lang = $ws.Range("D2").Value2
$prefix = $ws.Range("E2").Value2

$ws.Range("F1").Value = "Theme"

$themes = @(
    "SAP Belize Deep",
    "SAP High Contrast Black",
    "SAP High Contrast Black (SAP Belize)",
    "SAP High Contrast White",
    "SAP High Contrast White (SAP Belize)",
    "SAP Quartz Light",
    "SAP Belize"
)

for ($i = 6; $i -ge 1; $i--) {
    $row = 2 + $i
    $ws.Range("A2:D2").Copy()
    $ws.Range("A$row`:D$row").PasteSpecial(-4122)
    $ws.Range("E2").Copy()
    $ws.Range("F$row").PasteSpecial(-4122)
    $ws.Range("A2").Copy()
    $ws.Range("E$row").PasteSpecial(-4122)
}

# Now handle row 2 itself: E2 currently has the "last column" style (s=3);
# move that style to F2 (the new last column) and give E2 the regular body style (s=1)
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

for ($i = 0; $i -le 6; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value2 = $url
    $ws.Range("B$row").Value2 = $user
    $ws.Range("C$row").Value2 = $pass
    $ws.Range("D$row").Value2 = $lang
    $ws.Range("E$row").Value2 = $prefix
    $ws.Range("F$row").Value2 = $themes[$i]
}

$excel.CutCopyMode = 0
$ws.Columns.Item(6).ColumnWidth = 30
